# Populate the (previously empty) Accuracy/Loss columns I and J on the
# "StorageBox" worksheet with the NN analysis results, rows 11-46.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StorageBox")

# row, I value, J value
$data = @(
    @(11, 0.767, 0.487),
    @(12, 1.0, 0.001),
    @(13, 1.0, 0.001),
    @(14, 0.999, 0.007),
    @(15, 0.924, 0.351),
    @(16, 1.0, 0.0),
    @(17, 0.896, 0.406),
    @(18, 1.0, 0.001),
    @(19, 0.1, 2.303),
    @(20, 1.0, 0.001),
    @(21, 0.896, 0.287),
    @(22, 0.099, 2.303),
    @(23, 1.0, 0.001),
    @(24, 1.0, 0.001),
    @(25, 0.777, 0.409),
    @(26, 0.1, 2.303),
    @(27, 1.0, 0.004),
    @(28, 1.0, 0.001),
    @(29, 1.0, 0.001),
    @(30, 0.896, 0.152),
    @(31, 1.0, 0.002),
    @(32, 1.0, 0.001),
    @(33, 0.099, 0.203),
    @(34, 1.0, 0.002),
    @(35, 0.798, 0.415),
    @(36, 1.0, 0.001),
    @(37, 1.0, 0.017),
    @(38, 0.987, 0.051),
    @(39, 1.0, 0.003),
    @(40, 0.899, 0.307),
    @(41, 1.0, 0.003),
    @(42, 1.0, 0.001),
    @(43, 1.0, 0.003),
    @(44, 0.898, 0.151),
    @(45, 1.0, 0.003),
    @(46, 0.898, 0.148)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value2 = $row[1]
    $ws.Cells.Item($r, 10).Value2 = $row[2]
}

# Update the view state: scroll so row 3 is at top and select L11.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("L11").Select()
